# Scheduled runner update: refresh computed market-price / profit figures
# across several leve-profit sheets (static values, not formulas).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3463.7646
$ws.Range("I51").Value = 2440
$ws.Range("J51").Value = 4180.4
$ws.Range("K51").Value = 2440
$ws.Range("L51").Value = 4180.4
$ws.Range("M51").Value = -1956
$ws.Range("N51").Value = -5148.4
$ws.Range("H76").Value = 4276716
$ws.Range("I76").Value = 5294029
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 5294029
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -5293714
$ws.Range("N76").Value = -4630
$ws.Range("H79").Value = 4276716
$ws.Range("I79").Value = 5294029
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 5294029
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -5292937
$ws.Range("N79").Value = -6184
$ws.Range("H129").Value = 905.3333
$ws.Range("I129").Value = 329
$ws.Range("J129").Value = 1337.5834
$ws.Range("K129").Value = 987
$ws.Range("L129").Value = 4012.7502
$ws.Range("M129").Value = 4013
$ws.Range("N129").Value = -14012.7502
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2166.6667
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -1713
$ws.Range("N16").Value = -3074
$ws.Range("H31").Value = 5093.2285
$ws.Range("I31").Value = 1378.1052
$ws.Range("J31").Value = 9504.9375
$ws.Range("K31").Value = 1378.1052
$ws.Range("L31").Value = 9504.9375
$ws.Range("M31").Value = -1083.1052
$ws.Range("N31").Value = -10094.9375
$ws.Range("H34").Value = 5093.2285
$ws.Range("I34").Value = 1378.1052
$ws.Range("J34").Value = 9504.9375
$ws.Range("K34").Value = 1378.1052
$ws.Range("L34").Value = 9504.9375
$ws.Range("M34").Value = -1176.1052
$ws.Range("N34").Value = -9908.9375
$ws.Range("H58").Value = 2166.4
$ws.Range("I58").Value = 1368
$ws.Range("K58").Value = 1368
$ws.Range("M58").Value = -1165
$ws.Range("H113").Value = 2166.6667
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -6840
$ws.Range("H136").Value = 2166.4
$ws.Range("I136").Value = 1368
$ws.Range("K136").Value = 4104
$ws.Range("M136").Value = -1554
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 8300.4
$ws.Range("I58").Value = 1005
$ws.Range("J58").Value = 9111
$ws.Range("K58").Value = 3015
$ws.Range("L58").Value = 27333
$ws.Range("M58").Value = -2887
$ws.Range("N58").Value = -27589
$ws.Range("H75").Value = 1575.4117
$ws.Range("I75").Value = 735
$ws.Range("J75").Value = 1834
$ws.Range("K75").Value = 2205
$ws.Range("L75").Value = 5502
$ws.Range("M75").Value = -1207
$ws.Range("N75").Value = -7498
$ws.Range("H78").Value = 1575.4117
$ws.Range("I78").Value = 735
$ws.Range("J78").Value = 1834
$ws.Range("K78").Value = 6615
$ws.Range("L78").Value = 16506
$ws.Range("M78").Value = -1623
$ws.Range("N78").Value = -26490
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 27523
$ws.Range("J46").Value = 27523
$ws.Range("L46").Value = 27523
$ws.Range("N46").Value = -27835
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H70").Value = 4935.039
$ws.Range("J70").Value = 5158.385
$ws.Range("L70").Value = 5158.385
$ws.Range("N70").Value = -5698.385
$ws.Range("H73").Value = 4935.039
$ws.Range("J73").Value = 5158.385
$ws.Range("L73").Value = 5158.385
$ws.Range("N73").Value = -7030.385
$ws.Range("H80").Value = 2596.6667
$ws.Range("I80").Value = 2372.7273
$ws.Range("K80").Value = 2372.7273
$ws.Range("M80").Value = -1374.7273
$ws.Range("H83").Value = 2596.6667
$ws.Range("I83").Value = 2372.7273
$ws.Range("K83").Value = 11863.6365
$ws.Range("M83").Value = -6871.636500000001
$ws.Range("H113").Value = 1166.8572
$ws.Range("I113").Value = 822
$ws.Range("J113").Value = 1626.6666
$ws.Range("K113").Value = 822
$ws.Range("L113").Value = 1626.6666
$ws.Range("M113").Value = 1348
$ws.Range("N113").Value = -5966.6666
$ws.Range("H132").Value = 3304.639
$ws.Range("I132").Value = 3102.3076
$ws.Range("J132").Value = 3830.7
$ws.Range("K132").Value = 9306.9228
$ws.Range("L132").Value = 11492.1
$ws.Range("M132").Value = -6776.9228
$ws.Range("N132").Value = -16552.1
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2743.697
$ws.Range("I7").Value = 2122.4614
$ws.Range("J7").Value = 3147.5
$ws.Range("K7").Value = 2122.4614
$ws.Range("L7").Value = 3147.5
$ws.Range("M7").Value = -2010.4614
$ws.Range("N7").Value = -3371.5
$ws.Range("H61").Value = 1393.6154
$ws.Range("I61").Value = 1291.2
$ws.Range("J61").Value = 1735
$ws.Range("K61").Value = 1291.2
$ws.Range("L61").Value = 1735
$ws.Range("M61").Value = -1089.2
$ws.Range("N61").Value = -2139
$ws.Range("H113").Value = 1393.6154
$ws.Range("I113").Value = 1291.2
$ws.Range("J113").Value = 1735
$ws.Range("K113").Value = 1291.2
$ws.Range("L113").Value = 1735
$ws.Range("M113").Value = 878.8
$ws.Range("N113").Value = -6075
$ws.Range("H126").Value = 2743.697
$ws.Range("I126").Value = 2122.4614
$ws.Range("J126").Value = 3147.5
$ws.Range("K126").Value = 6367.3842
$ws.Range("L126").Value = 9442.5
$ws.Range("M126").Value = -3897.3842
$ws.Range("N126").Value = -14382.5
$ws.Range("H139").Value = 45571.668
$ws.Range("J139").Value = 45571.668
$ws.Range("L139").Value = 45571.668
$ws.Range("N139").Value = -55851.668
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 92190.17999999999
$ws.Range("I122").Value = 125824
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 377472
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -375022
$ws.Range("N122").Value = -12400
$ws.Range("H129").Value = 40000
$ws.Range("J129").Value = 40000
$ws.Range("L129").Value = 40000
$ws.Range("N129").Value = -50000
